$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$v = $ws.Range("A1").Value
Write-Host "A1=$v"
$v2 = $ws.Range("D36").Value
Write-Host "D36=$v2"
$v3 = $ws.Range("A35").Value
Write-Host "A35=$v3"
